$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "process" column (B) values of 3 are being relabeled as the new
# "final"/aggregate process marker "f" for DMUs A, B, C, D, E (rows 5, 9, 13, 17, 21)
$ws.Range("B5").Value = "f"
$ws.Range("B9").Value = "f"
$ws.Range("B13").Value = "f"
$ws.Range("B17").Value = "f"
$ws.Range("B21").Value = "f"

# Header J1 renamed from "y3" to "yf" to match the new "f" process naming
$ws.Range("J1").Value = "yf"

# Move the active selection to K6 (matches the saved workbook state)
$ws.Range("K6").Select()
